$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo: "Need do" -> "Need to do"
$ws.Range("A3").Value = "2. Need to do our best as individual and as a team."

# Update the selected cell to match the post-edit selection
$ws.Range("A13").Select()
